$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("2025-07-18 Friday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-07-19 Saturday", 2)
$null = $d.Content.Find.Execute("838×9=", $true, $true, $false, $false, $false, $true, 1, $false, "609×2=", 2)
$null = $d.Content.Find.Execute("769×5=", $true, $true, $false, $false, $false, $true, 1, $false, "850×6=", 2)
$null = $d.Content.Find.Execute("141×3=", $true, $true, $false, $false, $false, $true, 1, $false, "117×9=", 2)
$null = $d.Content.Find.Execute("879×7=", $true, $true, $false, $false, $false, $true, 1, $false, "411×7=", 2)
$null = $d.Content.Find.Execute("293×6=", $true, $true, $false, $false, $false, $true, 1, $false, "490×5=", 2)
$null = $d.Content.Find.Execute("635×2=", $true, $true, $false, $false, $false, $true, 1, $false, "546×9=", 2)
$null = $d.Content.Find.Execute("173×2=", $true, $true, $false, $false, $false, $true, 1, $false, "460×6=", 2)
$null = $d.Content.Find.Execute("768×5=", $true, $true, $false, $false, $false, $true, 1, $false, "868×5=", 2)
$null = $d.Content.Find.Execute("348×7=", $true, $true, $false, $false, $false, $true, 1, $false, "638×5=", 2)
$null = $d.Content.Find.Execute("834×4=", $true, $true, $false, $false, $false, $true, 1, $false, "520×3=", 2)
$null = $d.Content.Find.Execute("664×7=", $true, $true, $false, $false, $false, $true, 1, $false, "417×6=", 2)
$null = $d.Content.Find.Execute("167×6=", $true, $true, $false, $false, $false, $true, 1, $false, "625×2=", 2)
$null = $d.Content.Find.Execute("434×6=", $true, $true, $false, $false, $false, $true, 1, $false, "683×9=", 2)
$null = $d.Content.Find.Execute("350×3=", $true, $true, $false, $false, $false, $true, 1, $false, "612×5=", 2)
$null = $d.Content.Find.Execute("615×4=", $true, $true, $false, $false, $false, $true, 1, $false, "883×8=", 2)
$null = $d.Content.Find.Execute("252×4=", $true, $true, $false, $false, $false, $true, 1, $false, "248×5=", 2)
$null = $d.Content.Find.Execute("566×9=", $true, $true, $false, $false, $false, $true, 1, $false, "403×8=", 2)
$null = $d.Content.Find.Execute("683×2=", $true, $true, $false, $false, $false, $true, 1, $false, "776×7=", 2)
$null = $d.Content.Find.Execute("653×3=", $true, $true, $false, $false, $false, $true, 1, $false, "728×9=", 2)
$null = $d.Content.Find.Execute("946×4=", $true, $true, $false, $false, $false, $true, 1, $false, "642×2=", 2)
$null = $d.Content.Find.Execute("109×9=", $true, $true, $false, $false, $false, $true, 1, $false, "554×9=", 2)
$null = $d.Content.Find.Execute("746×2=", $true, $true, $false, $false, $false, $true, 1, $false, "192×5=", 2)
$null = $d.Content.Find.Execute("860×8=", $true, $true, $false, $false, $false, $true, 1, $false, "412×5=", 2)
$null = $d.Content.Find.Execute("943×2=", $true, $true, $false, $false, $false, $true, 1, $false, "283×4=", 2)
$null = $d.Content.Find.Execute("555×3=", $true, $true, $false, $false, $false, $true, 1, $false, "265×5=", 2)

$d.Save()
